# Converts an EMU (English Metric Unit) integer into the point value that,
# once PowerPoint's COM layer truncates it back to EMU (914400 EMU/in,
# 12700 EMU/pt), reproduces the exact original EMU integer.
function EMUToPt($emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape "TextBox 16" (id=17): colour every run accent6
# ---------------------------------------------------------------------
$sh17 = $s.Shapes.Item("TextBox 16")
$tr17 = $sh17.TextFrame.TextRange
$all17 = $tr17.Characters(1, $tr17.Length)
$all17.Font.Color.ObjectThemeColor = 10   # msoThemeColorAccent6

# ---------------------------------------------------------------------
# Shape "TextBox 21" (id=22): colour every run red (FF0000)
# ---------------------------------------------------------------------
$sh22 = $s.Shapes.Item("TextBox 21")
$tr22 = $sh22.TextFrame.TextRange
$all22 = $tr22.Characters(1, $tr22.Length)
$all22.Font.Color.RGB = 255   # RGB(255,0,0) -> &H0000FF (BGR packed) -> FF0000 in OOXML

# ---------------------------------------------------------------------
# Shape "Straight Arrow Connector 28" (id=29): dash the line (sysDot)
# Setting DashStyle alone appends <a:prstDash> at the end of <a:ln>; also
# touching the (unchanged) arrowhead styles forces the line properties to
# be re-emitted in schema order, putting <a:prstDash> before <a:headEnd>/
# <a:tailEnd> as PowerPoint itself would.
# ---------------------------------------------------------------------
$sh29 = $s.Shapes.Item("Straight Arrow Connector 28")
$sh29.Line.DashStyle = 2          # msoLineSquareDot -> prstDash val="sysDot" (engine mapping)
$sh29.Line.BeginArrowheadStyle = $sh29.Line.BeginArrowheadStyle
$sh29.Line.EndArrowheadStyle = $sh29.Line.EndArrowheadStyle

# ---------------------------------------------------------------------
# Shape "TextBox 33" (id=34): grow the box, prefix first run with
# "(optional) ", colour every run accent1
# ---------------------------------------------------------------------
$sh34 = $s.Shapes.Item("TextBox 33")
$sh34.Height = EMUToPt 1754326

$tr34 = $sh34.TextFrame.TextRange
$firstChar34 = $tr34.Characters(1, 1)
$firstChar34.Text = "(optional) C"

$tr34b = $sh34.TextFrame.TextRange
$all34 = $tr34b.Characters(1, $tr34b.Length)
$all34.Font.Color.ObjectThemeColor = 5   # msoThemeColorAccent1

# ---------------------------------------------------------------------
# Shape "TextBox 35" (id=36): reposition only
# ---------------------------------------------------------------------
$sh36 = $s.Shapes.Item("TextBox 35")
$sh36.Left = EMUToPt 1997130
$sh36.Top = EMUToPt 3629129

# ---------------------------------------------------------------------
# Shape "TextBox 36" (id=37): reposition, widen, fix typo in the text
# ---------------------------------------------------------------------
$sh37 = $s.Shapes.Item("TextBox 36")
$sh37.Left = EMUToPt 5086809
$sh37.Top = EMUToPt 3609536
$sh37.Width = EMUToPt 2771977

$tr37 = $sh37.TextFrame.TextRange
$tr37.Text = "NoAxDsimVSNoAxDpred"
